# edit.ps1 - PowerPoint COM-interop script
#
# Applies the two substantive edits captured by the commit diff
# (the repeated 2/4/2025 -> 2/5/2025 hunks are just PowerPoint's
# auto-updating "datetimeFigureOut" field re-caching itself because the
# deck was re-saved a day later; that is not a deliberate content edit,
# so it is intentionally left alone here):
#
#   1. Slide 1, "Title 1" (ctrTitle) placeholder: reflowed from a tall
#      box near the top of the slide to a short, word-wrapped,
#      auto-fit-to-text box positioned further down.
#   2. Slide 2, "Content Placeholder 2": appended a "Products:"/"-"/"-"
#      "Services:"/"-"/"-" scaffold after the existing answer paragraph.

$p = $ppt.ActivePresentation

# --- 1. Slide 1 title placeholder resize/reposition -----------------------
$slide1 = $p.Slides.Item(1)
$title = $slide1.Shapes.Item(3)   # "Title 1" (ph type="ctrTitle")

# Flip autofit behaviour first (noAutofit -> spAutoFit, word wrap on) so the
# later explicit Height assignment below is not clobbered by an autosize
# recalculation happening afterwards.
$title.TextFrame.WordWrap = 1
$title.TextFrame.AutoSize = 1

# Values are expressed in points (English Metric Units / 12700) to land the
# shape at the exact target offsets/extents used by the canonical deck.
$title.Left = 60.89827771653544
$title.Top = 229.67111236220472
$title.Width = 284.7963192125984
$title.Height = 68.53788401574803

# --- 2. Slide 2 content placeholder: append Products/Services scaffold ----
$slide2 = $p.Slides.Item(2)
$content = $slide2.Shapes.Item(2)  # "Content Placeholder 2"

$tr = $content.TextFrame.TextRange
$tr.InsertAfter("`rProducts:`r- `r- `rServices:`r-`r-")
